$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet tracks fixtures for "Thu Oct 23"; this update records results for
# two more finished matches and reorders/refreshes the list of fixtures
# (newly-finished matches bubble up, already-finished ones shift down, and
# two brand-new rows are appended before the summary block).
#
# Insert two new rows before row 13 so the trailing summary-formula rows
# (old rows 13-15) shift down to rows 15-17; Excel auto-adjusts the
# formulas in the process (K13->K15, K14->K16, K15->K17, and the
# K13/K15 references inside them).
$ws.Rows("13:14").Insert()

# Row 2
$ws.Range("A2").Value = "Thu Oct 23"
$ws.Range("B2").Value = "FC Krasnodar ✓ - FC Sochi: 3:0"
$ws.Range("C2").Value = 3.56
$ws.Range("D2").Value = "FC Krasnodar"
$ws.Range("E2").Value = 4.5
$c = $ws.Range("F2")
$c.NumberFormat = "@"
$c.Value = "72%"
$c.Style = "Normal"
$ws.Range("G2").Value = "✓"
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = $True

# Row 3
$ws.Range("A3").Value = "Thu Oct 23"
$ws.Range("B3").Value = "Al-Najma SC - Al-Ahli SFC ✓: 0:1"
$ws.Range("C3").Value = 2.12
$ws.Range("D3").Value = "Al-Ahli SFC"
$ws.Range("E3").Value = 3.5
$c = $ws.Range("F3")
$c.NumberFormat = "@"
$c.Value = "66%"
$c.Style = "Normal"
$ws.Range("G3").Value = "✓"
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = $True

# Row 4
$ws.Range("A4").Value = "Thu Oct 23"
$ws.Range("B4").Value = "FC Tallinn - FC Nomme United ✓: 1:4"
$ws.Range("C4").Value = 3.09
$ws.Range("D4").Value = "FC Nomme United"
$ws.Range("E4").Value = 4.5
$c = $ws.Range("F4")
$c.NumberFormat = "@"
$c.Value = "65%"
$c.Style = "Normal"
$ws.Range("G4").Value = "✓"
$ws.Range("H4").Value = 5
$ws.Range("I4").Value = $False

# Row 5
$ws.Range("A5").Value = "Thu Oct 23"
$ws.Range("B5").Value = "CR Flamengo ✓ - Racing Club: 1:0"
$ws.Range("C5").Value = 1.59
$ws.Range("D5").Value = "CR Flamengo"
$ws.Range("E5").Value = 2.5
$c = $ws.Range("F5")
$c.NumberFormat = "@"
$c.Value = "65%"
$c.Style = "Normal"
$ws.Range("G5").Value = "✓"
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = $True

# Row 6
$ws.Range("A6").Value = "Thu Oct 23"
$ws.Range("B6").Value = "FC Rànger's  - FC Pas de la Casa: 13:30"
$ws.Range("C6").Value = 2.14
$ws.Range("D6").Value = "FC Rànger's"
$ws.Range("E6").Value = 3.5
$c = $ws.Range("F6")
$c.NumberFormat = "@"
$c.Value = "64%"
$c.Style = "Normal"
$ws.Range("G6").ClearContents()
$ws.Range("H6").Value = 43
$ws.Range("I6").Value = $False

# Row 7
$ws.Range("A7").Value = "Thu Oct 23"
$ws.Range("B7").Value = "Maccabi Tel Aviv - FC Midtjylland ✓: 0:3"
$ws.Range("C7").Value = 2.63
$ws.Range("D7").Value = "FC Midtjylland"
$ws.Range("E7").Value = 3.5
$c = $ws.Range("F7")
$c.NumberFormat = "@"
$c.Value = "64%"
$c.Style = "Normal"
$ws.Range("G7").Value = "✓"
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = $True

# Row 8
$ws.Range("A8").Value = "Thu Oct 23"
$ws.Range("B8").Value = "Feyenoord Rotterdam ✓ - Panathinaikos FC: 3:1"
$ws.Range("C8").Value = 1.95
$ws.Range("D8").Value = "Feyenoord Rotterdam"
$ws.Range("E8").Value = 2.5
$c = $ws.Range("F8")
$c.NumberFormat = "@"
$c.Value = "63%"
$c.Style = "Normal"
$ws.Range("G8").Value = "✓"
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = $False

# Row 9
$ws.Range("A9").Value = "Thu Oct 23"
$ws.Range("B9").Value = "Atlético Nacional ✓ - Once Caldas: 2:0"
$ws.Range("C9").Value = 1.75
$ws.Range("D9").Value = "Atlético Nacional"
$ws.Range("E9").Value = 2.5
$c = $ws.Range("F9")
$c.NumberFormat = "@"
$c.Value = "61%"
$c.Style = "Normal"
$ws.Range("G9").Value = "✓"
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = $True

# Row 10
$ws.Range("A10").Value = "Thu Oct 23"
$ws.Range("B10").Value = "Club Tijuana - Deportivo Toluca : 0:0"
$ws.Range("C10").Value = 3.21
$ws.Range("D10").Value = "Deportivo Toluca"
$ws.Range("E10").Value = 4.5
$c = $ws.Range("F10")
$c.NumberFormat = "@"
$c.Value = "61%"
$c.Style = "Normal"
$ws.Range("G10").ClearContents()
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = $True

# Row 11
$ws.Range("A11").Value = "Thu Oct 23"
$ws.Range("B11").Value = "Xorazm Urganch - Pakhtakor Tashkent ✓: 0:1"
$ws.Range("C11").Value = 1.95
$ws.Range("D11").Value = "Pakhtakor Tashkent"
$ws.Range("E11").Value = 2.5
$c = $ws.Range("F11")
$c.NumberFormat = "@"
$c.Value = "60%"
$c.Style = "Normal"
$ws.Range("G11").Value = "✓"
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = $True

# Row 12
$ws.Range("A12").Value = "Thu Oct 23"
$ws.Range("B12").Value = "Shamrock Rovers - NK Celje ✓: 0:2"
$ws.Range("C12").Value = 3.28
$ws.Range("D12").Value = "NK Celje"
$ws.Range("E12").Value = 4.5
$c = $ws.Range("F12")
$c.NumberFormat = "@"
$c.Value = "58%"
$c.Style = "Normal"
$ws.Range("G12").Value = "✓"
$ws.Range("H12").Value = 2
$ws.Range("I12").Value = $True

# Row 13
$ws.Range("A13").Value = "Thu Oct 23"
$ws.Range("B13").Value = "AEK Athens ✓ - Aberdeen FC: 6:0"
$ws.Range("C13").Value = 2.34
$ws.Range("D13").Value = "AEK Athens"
$ws.Range("E13").Value = 3.5
$c = $ws.Range("F13")
$c.NumberFormat = "@"
$c.Value = "56%"
$c.Style = "Normal"
$ws.Range("G13").Value = "✓"
$ws.Range("H13").Value = 6
$ws.Range("I13").Value = $False

# Row 14
$ws.Range("A14").Value = "Thu Oct 23"
$ws.Range("B14").Value = "Celtic FC ✓ - SK Sturm Graz: 2:1"
$ws.Range("C14").Value = 2.5
$ws.Range("D14").Value = "Celtic FC"
$ws.Range("E14").Value = 3.5
$c = $ws.Range("F14")
$c.NumberFormat = "@"
$c.Value = "54%"
$c.Style = "Normal"
$ws.Range("G14").Value = "✓"
$ws.Range("H14").Value = 3
$ws.Range("I14").Value = $True

